# Insert a new data row at row 270 (shifting the existing rows 270-308 down
# to 271-309), then populate the new row. The "fixed" columns (market /
# product descriptors) for the new row are identical to the row that used to
# be at 270 (now at 271), so they are copied across; only the date, volume,
# price and origin columns carry genuinely new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 270:308 down to 271:309, leaving row 270 empty.
$ws.Rows.Item(270).EntireRow.Insert()

# Columns that stay identical to the (now shifted) old row 270 / new row 271.
$fixedCols = @("A", "B", "C", "E", "F", "G", "H", "I", "J", "K", "L", "Q", "T")
foreach ($col in $fixedCols) {
    $ws.Range($col + "270").Value2 = $ws.Range($col + "271").Value2
}

# New values for the inserted row.
$ws.Range("D270").Value2 = 44984
$ws.Range("M270").Value2 = 250
$ws.Range("N270").Value2 = 3000
$ws.Range("O270").Value2 = 3000
$ws.Range("P270").Value2 = 3000
$ws.Range("R270").Value2 = "Perú"
$ws.Range("S270").Value2 = 1500
